# Regenerate merged AHB files
# - Rename the "_old" / "_new" header suffixes to "_FV2210" / "_FV2304"
# - Turn the A1:U58 range into a real Excel Table ("Table1")
# - Freeze the header row (split/freeze below row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:J1 = "_old" -> "_FV2210", L1:U1 = "_new" -> "_FV2304") ---
$ws.Range("A1").Value = "Segmentname_FV2210"
$ws.Range("B1").Value = "Segmentgruppe_FV2210"
$ws.Range("C1").Value = "Segment_FV2210"
$ws.Range("D1").Value = "Datenelement_FV2210"
$ws.Range("E1").Value = "Segment ID_FV2210"
$ws.Range("F1").Value = "Code_FV2210"
$ws.Range("G1").Value = "Qualifier_FV2210"
$ws.Range("H1").Value = "Beschreibung_FV2210"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2210"
$ws.Range("J1").Value = "Bedingung_FV2210"

$ws.Range("L1").Value = "Segmentname_FV2304"
$ws.Range("M1").Value = "Segmentgruppe_FV2304"
$ws.Range("N1").Value = "Segment_FV2304"
$ws.Range("O1").Value = "Datenelement_FV2304"
$ws.Range("P1").Value = "Segment ID_FV2304"
$ws.Range("Q1").Value = "Code_FV2304"
$ws.Range("R1").Value = "Qualifier_FV2304"
$ws.Range("S1").Value = "Beschreibung_FV2304"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("U1").Value = "Bedingung_FV2304"

# --- 2. Convert the used range into an Excel Table with header row ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), $false, 1)
$tbl.Name = "Table1"

# --- 3. Freeze panes below the header row ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
